$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = $origStyle
}

$ws.Range("D2").Value = "25.880.51"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "1.584.11"
$ws.Range("E3").Value = "  -2.22%  "
$ws.Range("E4").Value = "  -0.26%  "
Set-TextValue "D5" "210.01"
$ws.Range("E5").Value = "  -1.22%  "
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("E7").Value = "  -2.52%  "
$ws.Range("E8").Value = "  -0.64%  "
$ws.Range("E9").Value = "  -1.19%  "
Set-TextValue "D10" "18.09"
$ws.Range("E10").Value = "  -1.18%  "
Set-TextValue "D11" "0.0790"
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").Value = "1.802.88"
$ws.Range("E12").Value = "  -2.31%  "
$ws.Range("D13").Value = "1.577.07"
$ws.Range("E13").Value = "  -2.63%  "
$ws.Range("E14").Value = "  -2.73%  "
$ws.Range("E15").Value = "  -2.46%  "
$ws.Range("D16").Value = "25.845.46"
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("D17").Value = "0.0₃0726"
$ws.Range("E17").Value = "  -1.62%  "
Set-TextValue "D18" "59.75"
$ws.Range("E18").Value = "  -3.18%  "
$ws.Range("E19").Value = "  -0.20%  "
Set-TextValue "D20" "191.90"
$ws.Range("E20").Value = "  +0.05%  "
Set-TextValue "D21" "4.18"
$ws.Range("E21").Value = "  -1.66%  "
$ws.Range("E22").Value = "  -1.63%  "
$ws.Range("E23").Value = "  -1.36%  "
$ws.Range("E24").Value = "  +0.10%  "
Set-TextValue "D25" "141.54"
$ws.Range("E25").Value = "  -1.96%  "
$ws.Range("E26").Value = "  -0.24%  "
$ws.Range("E27").Value = "  -1.43%  "
Set-TextValue "D28" "15.07"
$ws.Range("E28").Value = "  -0.90%  "
$ws.Range("E30").Value = "  -5.50%  "
Set-TextValue "D31" "0.0471"
$ws.Range("E31").Value = "  -1.45%  "
$ws.Range("E32").Value = "  +0.25%  "
$ws.Range("E33").Value = "  -2.76%  "
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("E35").Value = "  -2.43%  "
$ws.Range("D36").Value = "1.096.96"
$ws.Range("E36").Value = "  -2.72%  "
$ws.Range("E37").Value = "  -0.31%  "
$ws.Range("E38").Value = "  -2.37%  "
$ws.Range("E39").Value = "  -2.00%  "
Set-TextValue "D40" "0.501"
$ws.Range("E40").Value = "  -3.48%  "
$ws.Range("E41").Value = "  -8.02%  "
Set-TextValue "D42" "0.812"
$ws.Range("E42").Value = "  +7.56%  "
Set-TextValue "D43" "93.81"
$ws.Range("E43").Value = "  -4.13%  "
Set-TextValue "D44" "5.14"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").Value = "1.715.75"
$ws.Range("E45").Value = "  -2.33%  "
$ws.Range("E46").Value = "  -1.15%  "
$ws.Range("E47").Value = "  -0.26%  "
Set-TextValue "D48" "53.12"
$ws.Range("E48").Value = "  -1.74%  "
$ws.Range("E49").Value = "  -1.38%  "
$ws.Range("E50").Value = "  -0.82%  "
$ws.Range("E51").Value = "  -0.34%  "
